$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings are not
# auto-converted to numbers; style is restored to Normal afterwards so the
# cells end up with no explicit style, matching the original workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.180.88'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '2.570.58'
$ws.Range("E3").Value = '  +1.73%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '583.18'
$ws.Range("E5").Value = '  +2.68%  '
$ws.Range("D6").Value = '148.47'
$ws.Range("E6").Value = '  +1.24%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("E9").Value = '  +2.88%  '
$ws.Range("D10").Value = '5.59'
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").Value = '27.85'
$ws.Range("E13").Value = '  +0.89%  '
$ws.Range("D14").Value = '3.033.96'
$ws.Range("E14").Value = '  +1.73%  '
$ws.Range("D15").Value = '63.119.52'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("D17").Value = '2.578.80'
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("D18").Value = '11.45'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").Value = '341.71'
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("E21").Value = '  +1.82%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").Value = '65.99'
$ws.Range("E23").Value = '  +1.06%  '
$ws.Range("D24").Value = '2.686.19'
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("E25").Value = '  +3.65%  '
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("D27").Value = '8.24'
$ws.Range("E27").Value = '  +14.55%  '
$ws.Range("D28").Value = '8.53'
$ws.Range("E28").Value = '  +2.43%  '
$ws.Range("B29").Value = 'SuiNetwork'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D29").Value = '1.49'
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").Value = '1.96'
$ws.Range("E31").Value = '  +6.14%  '
$ws.Range("D32").Value = '0.0₃0828'
$ws.Range("E32").Value = '  +2.52%  '
$ws.Range("D33").Value = '177.22'
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").Value = '436.09'
$ws.Range("E34").Value = '  +5.56%  '
$ws.Range("D35").Value = '1.60'
$ws.Range("E35").Value = '  +1.24%  '
$ws.Range("D36").Value = '0.407'
$ws.Range("E36").Value = '  +2.44%  '
$ws.Range("D37").Value = '19.31'
$ws.Range("E37").Value = '  +2.75%  '
$ws.Range("D38").Value = '4.54'
$ws.Range("E38").Value = '  +3.50%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  +0.74%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = '152.21'
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("D44").Value = '21.39'
$ws.Range("E44").Value = '  +3.83%  '
$ws.Range("D45").Value = '0.0555'
$ws.Range("E45").Value = '  +6.70%  '
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").Value = '0.0977'
$ws.Range("E47").Value = '  +1.23%  '
$ws.Range("D48").Value = '0.0242'
$ws.Range("E48").Value = '  +2.27%  '
$ws.Range("D49").Value = '18.49'
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("D51").Value = '11.37'
$ws.Range("E51").Value = '  -0.25%  '

$ws.Range("D2:D51").Style = "Normal"
